# Adds the 2020 data column (Q) to the 3.3.4 Hepatitis B incidence sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2020 values (column Q) for each data row, keyed by row number.
$q2020 = @{
    3 = 2020
    4 = 1.9148453093736542
    5 = 1.7453236044300597
    6 = 2.0818900906859255
    7 = 1.658050942694075
    8 = 1.4467487937731931
    9 = 1.8774124750304142
    10 = 0.96024351775610284
    11 = 0.63595936855594293
    12 = 1.2888424905592288
    13 = 1.6032353288937073
    14 = 2.4146715443031859
    15 = 0.79837132250209564
    16 = 1.3751327862596732
    17 = 0.67516929870164943
    18 = 2.1012817818869509
    19 = 1.5943738893736428
    20 = 1.5765365498500856
    21 = 1.6126194804433236
    22 = 0.37150276583809166
    23 = 0
    24 = 0.75125835774923
    25 = 2.8942542850468351
    26 = 2.72898263527357
    27 = 3.0545792215303034
    28 = 3.9473869708034344
    29 = 3.6031203021816895
    30 = 4.2520923837938582
    31 = 0
    32 = 0
    33 = 0
}

foreach ($r in 3..33) {
    # Copy column P's formatting (number format/font/border) into column Q
    # for this row, then overwrite the pasted value with the 2020 figure.
    $ws.Range("P$r").Copy()
    $ws.Range("Q$r").PasteSpecial(-4122)
    $ws.Range("Q$r").Value = $q2020[$r]
}

$excel.CutCopyMode = $false

# Match the author's final cursor/selection position.
$ws.Range("T1").Select() | Out-Null
